$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update A26:A37 from 20 to 17
$ws.Range("A26:A37").Value = 17

# Update B26 from 18 to 0 (B27:B37 are formulas referencing previous row, they recalc automatically)
$ws.Range("B26").Value = 0

# Update the active selection on the sheet to B27
$ws.Range("B27").Select()
